# Applies the OOXML changes described by the commit:
#   1. Update the (fixed) date placeholder text "11/11/2020" -> "11/22/2020"
#      on the slide master and every slide layout.
#   2. Move the full-bleed background picture on slide 1 up so it sits
#      flush with the top of the slide (y offset 1673 EMU -> 0).
#   3. Fix wording/typos in the "Using Labels in Control Flow Statements"
#      slide (slide 8) body text box.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/11/2020") {
                $tr.Text = "11/22/2020"
            }
        }
    }
}

# --- 1. Date placeholder on the master and every layout ---------------
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2. Reposition the background picture on slide 1 -------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "Picture 4") {
        $shp.Left = 0
        $shp.Top = 0
    }
}

# --- 3. Fix the label-description text on slide 8 ----------------------
$slide8 = $p.Slides.Item(8)
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "A label is a simple word*") {
        $shp.TextFrame.TextRange.Text = "A label in a simple word is an identifier followed by a colon. It is applied to a statement or a block of code. Labels are mostly used when continue and break statements need to jump to certain block of code or to certain iterations."
    }
}
